# Append a fresh price-history snapshot to the holdings sheet.
#
# The previous run's snapshot (rows 2-9, all dated 2021-07-31) is replaced
# with a new snapshot dated 2021-08-08 (serial 44416.96780092592). The
# currency ordering also changes: BTC3S now sits right after BTC, a brand
# new BTC3L row is inserted after it, and a brand new MATIC row is appended
# at the bottom (row 10) so nothing is lost.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDate = 44416.96780092592

# Create row 10 first by cloning the formatting of row 9 (same column-A
# date style, s="2") before we start writing values into it.
$ws.Cells.Item(9, 1).Copy($ws.Cells.Item(10, 1))

# currency, balance, available, holds, price, dollar_value for every data
# row (2..10). "price" is $null where the column should stay blank.
$rows = @(
    @{ Row = 2;  Currency = "USDT";  Balance = 1326.50294401; Available = 795.32972885;  Holds = 531.17321516; Price = $null;    Dollar = 1326.5 },
    @{ Row = 3;  Currency = "USDC";  Balance = 1223.29700881; Available = 1223.29700881; Holds = 0;            Price = $null;    Dollar = 1223.3 },
    @{ Row = 4;  Currency = "BTC";   Balance = 0.02448253;    Available = 0.02448253;    Holds = 0;            Price = 44344.6;  Dollar = 1085.67 },
    @{ Row = 5;  Currency = "BTC3S"; Balance = 116.1203;      Available = 1.1076;        Holds = 115.0127;     Price = $null;    Dollar = 116.12 },
    @{ Row = 6;  Currency = "BTC3L"; Balance = 114.3699;      Available = 3.8651;        Holds = 110.5048;     Price = $null;    Dollar = 114.37 },
    @{ Row = 7;  Currency = "ATOM";  Balance = 0.2366;        Available = 0.2366;        Holds = 0;            Price = 13.402;   Dollar = 3.17 },
    @{ Row = 8;  Currency = "ALGO";  Balance = 0;              Available = 0;            Holds = 0;            Price = 0.8378;   Dollar = 0 },
    @{ Row = 9;  Currency = "ETH";   Balance = 0;              Available = 0;            Holds = 0;            Price = 3054.33;  Dollar = 0 },
    @{ Row = 10; Currency = "MATIC"; Balance = 0;              Available = 0;            Holds = 0;            Price = $null;    Dollar = 0 }
)

foreach ($row in $rows) {
    $r = $row.Row

    $ws.Cells.Item($r, 1).Value = $newDate          # A: date
    $ws.Cells.Item($r, 2).Value = $row.Currency     # B: currency
    $ws.Cells.Item($r, 3).Value = "trade"           # C: act_name
    $ws.Cells.Item($r, 4).Value = $row.Balance      # D: balance
    $ws.Cells.Item($r, 5).Value = $row.Available    # E: available
    $ws.Cells.Item($r, 6).Value = $row.Holds        # F: holds

    if ($null -eq $row.Price) {
        $ws.Cells.Item($r, 7).ClearContents()       # G: price (blank)
    } else {
        $ws.Cells.Item($r, 7).Value = $row.Price    # G: price
    }

    $ws.Cells.Item($r, 8).Value = $row.Dollar       # H: dollar_value
}
